# changed currency statement format
# Removes the "Döviz Borç / Döviz Alacak / Döviz Bakiye" (foreign-currency)
# header columns from row 4 (F4:H4) of the partner currency statement,
# re-centers the M2:N2 merged cell vertically, and moves the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Döviz Borç" / "Döviz Alacak" / "Döviz Bakiye" header text —
# these columns are no longer shown, the remaining headers (B/A, Kur,
# TL Borç, TL Alacak, TL Bakiye, B/A) shift to fill I4:N4 as before.
$ws.Range("F4:H4").ClearContents()

# M2:N2 (merged cell to the right of "Tel :") switches from top-aligned
# to vertically centered, matching the other header cells on row 2.
$ws.Range("M2:N2").VerticalAlignment = -4108

# Active cell / selection moved.
$ws.Range("H14").Select() | Out-Null
